# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# The two worker rows in the statement (rows 16 and 17) are updated: the
# previous account-statement entries are removed and replaced, which in
# practice swaps which worker's data sits in row 16 vs row 17.
#
# Row 16 now holds EDUAR ENRIQUE ARRIETA HERRERA (CC 10887028)
# Row 17 now holds EDER LUIS BUELVAS HERNANDEZ (CC 1047387735)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C16").Value = "10887028"
$ws.Range("D16").Value = "EDUAR ENRIQUE ARRIETA HERRERA"

$ws.Range("C17").Value = "1047387735"
$ws.Range("D17").Value = "EDER LUIS BUELVAS HERNANDEZ"
